$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.09181472831639451, 1.838119433117105, 17.95916271619692, 4.237825234267799, 4.278989066040988, 51)
    3  = @(0.00007413020273047021, 2.091730062774077, 18.37008874755701, 4.286034151468816, 4.329548303496322, 50)
    4  = @(0.07697517700982162, 1.786816454232815, 10.61989574786886, 3.258818152009844, 3.291670541969232, 49)
    5  = @(0.005243764112333811, 1.747033754236641, 11.76101659567932, 3.429433859353365, 3.465721123983597, 48)
    6  = @(0.05952959428276366, 2.303747599184862, 19.73205291551434, 4.442077544968608, 4.489698183225748, 47)
    7  = @(0.0295772053651569, 2.295364381528962, 15.73071993611709, 3.96619716304133, 4.00991237297396, 46)
    8  = @(0.08692572624757923, 2.233491525293952, 14.67667969733654, 3.83101549165969, 3.873307732227242, 45)
    9  = @(0.04420571028408105, 2.418601878655092, 15.97758211506545, 3.997196782129377, 4.04316141182984, 44)
    10 = @(0.1200444304270716, 2.442402265346571, 16.98616212122741, 4.121427194701782, 4.168433851739839, 43)
    11 = @(0.0538547395322746, 2.439448375966829, 16.12360189523408, 4.015420512877085, 4.063728510277974, 42)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
